$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extraneous row (old row 13: unlabeled row holding the
# docente name) so everything below shifts up by one row, matching the
# new A1:C21 layout (was A1:C22).
$ws.Rows("13:13").Delete()

# After the shift, a handful of B/C cells need their text corrected so
# that label (column A) and value (columns B/C) pairings match the
# target content.
$ws.Range("B10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2016" looks like a date, and a plain .Value assignment would
# get auto-converted to a date serial (changing both the stored type
# and the cell style). Route it through a text formula + paste-special
# values so it lands back in the sheet as plain text, preserving the
# existing cell style.
$ws.Range("B15").Formula = "=""01/01/2016"""
$ws.Range("C15").Formula = "=""01/01/2016"""
$ws.Range("B15:C15").Copy() | Out-Null
$ws.Range("B15:C15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"

$ws.Range("B19").Value = "A avaliação da disciplina será feita por meio de avaliação escrita e avaliação de atividades em grupo (seminários)"
$ws.Range("C19").Value = "A avaliação da disciplina será feita por meio de avaliação escrita e avaliação de atividades em grupo (seminários)"

$ws.Range("B20").Value = "A média será composta da seguinte forma:`nNota do semestre (M) = 0,6P+0,4S. Sendo uma prova na final do semestre (P) e o seminário (S) apresentado pelos alunos."
$ws.Range("C20").Value = "A média será composta da seguinte forma:`nNota do semestre (M) = 0,6P+0,4S. Sendo uma prova na final do semestre (P) e o seminário (S) apresentado pelos alunos."

$ws.Range("B21").Value = "A Nota de Recuperação (NR) será dada pela média aritmética entre a Nota do Semestre (M) e a Prova de Recuperação (PR), sendo aprovado o aluno que obtiver NR maior ou igual a cinco."
$ws.Range("C21").Value = "A Nota de Recuperação (NR) será dada pela média aritmética entre a Nota do Semestre (M) e a Prova de Recuperação (PR), sendo aprovado o aluno que obtiver NR maior ou igual a cinco."
